# This script edits /tmp/work/before.docx per the target diff:
#  1. Paragraph 1: merge split runs into a single run (no text change).
#  2. Paragraph 3 ("its special ..."): split into four paragraphs, trimming
#     text to stop before " (every ..." and marking "every" with
#     proofErr gramStart/gramEnd, and merging what used to be paragraph 4
#     ("I want some more features ...") in as a run-merged paragraph.
#  3. Last "kitchen staff" paragraph: merge split runs into a single run
#     (keeping the existing proofErr pair at the start intact).

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($range, [string]$bodyXml) {
    $xml = $pkgHeader + '<w:body>' + $bodyXml + '</w:body>' + $pkgFooter
    $range.InsertXML($xml)
}

# --- Edit 1: merge the runs of the first paragraph ---------------------
$r1 = $d.Content
$r1.Find.Execute("I wan", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$p1 = $r1.Paragraphs(1)
Set-ParagraphXml $p1.Range '<w:p><w:r><w:t xml:space="preserve">I want do to for me a simple application simple interfaces and connect with database. </w:t></w:r></w:p>'

# --- Edit 2: split the "its special ..." paragraph into four -----------
$r2 = $d.Content
$r2.Find.Execute("its special for our country", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$p3 = $r2.Paragraphs(1)
$newBody = '<w:p><w:r><w:t>its special for our country the scenario is the customer and the waiter</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>every</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> on have separate interface)</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t xml:space="preserve"> can order from the menu then the order send to the kitchen staff directly.</w:t></w:r></w:p>' + `
           '<w:p><w:r><w:t>I want some more features add to the app more than regular food app that the customer when arrived to the restaurant can push button then a notification sent to the kitchen.</w:t></w:r></w:p>'

# This paragraph's Range must also consume the following paragraph (the
# old "I want some more features ..." one) since it is being folded in
# and removed as a standalone paragraph.
$nextPara = $p3.Next()
$mergedRange = $d.Range($p3.Range.Start, $nextPara.Range.End)
Set-ParagraphXml $mergedRange $newBody

# --- Edit 3: merge the runs of the final "kitchen staff" paragraph -----
$r3 = $d.Content
$r3.Find.Execute("kitchen staff have button", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0) | Out-Null
$p5 = $r3.Paragraphs(1)
Set-ParagraphXml $p5.Range '<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>also  the</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> kitchen staff have button that  is order is ready then a notification send to the waiter to take the order to the customer.</w:t></w:r></w:p>'
